$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited cells keep their original text storage (many look numeric,
# e.g. "310.64", and Excel would silently coerce them to the Number type on
# assignment otherwise) by forcing a Text number format before writing values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.317.10"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.305.08"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.64"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.23"
$ws.Range("E6").Value = "  +6.02%  "
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +7.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.56"
$ws.Range("E10").Value = "  +3.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.33"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("E14").Value = "  +3.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.661.31"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.06"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.305.49"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.206.40"
$ws.Range("E19").Value = "  +2.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.17"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0925"
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("E22").Value = "  +2.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.06"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "242.25"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.61"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.92"
$ws.Range("E28").Value = "  +5.71%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  +8.01%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.89"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.55"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.27"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.27"
$ws.Range("E35").Value = "  +4.09%  "
$ws.Range("E36").Value = "  +6.52%  "
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.45"
$ws.Range("E41").Value = "  +7.69%  "
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("E43").Value = "  +14.99%  "
$ws.Range("E44").Value = "  +3.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.984.95"
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.85"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.02"
$ws.Range("E47").Value = "  +3.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.00"
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.81"
$ws.Range("E49").Value = "  +3.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.92"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("E51").Value = "  +7.90%  "
